$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H2").Value = 300.25
$ws.Range("I2").Value = 162.84616
$ws.Range("J2").Value = 895.6667
$ws.Range("K2").Value = 162.84616
$ws.Range("L2").Value = 895.6667
$ws.Range("M2").Value = -49.84616
$ws.Range("N2").Value = -1121.6667
$ws.Range("H38").Value = 314.375
$ws.Range("I38").Value = 314.375
$ws.Range("K38").Value = 943.125
$ws.Range("M38").Value = -571.125
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()
$ws.Range("H55").Value = 1194.7858
$ws.Range("I55").Value = 2439
$ws.Range("J55").Value = 503.55554
$ws.Range("K55").Value = 2439
$ws.Range("L55").Value = 503.55554
$ws.Range("M55").Value = -2225
$ws.Range("N55").Value = -931.5555400000001
$ws.Range("H88").Value = 1600.6154
$ws.Range("I88").Value = 2479
$ws.Range("J88").Value = 1051.625
$ws.Range("K88").Value = 2479
$ws.Range("L88").Value = 1051.625
$ws.Range("M88").Value = -2073
$ws.Range("N88").Value = -1863.625
$ws.Range("H91").Value = 1600.6154
$ws.Range("I91").Value = 2479
$ws.Range("J91").Value = 1051.625
$ws.Range("K91").Value = 2479
$ws.Range("L91").Value = 1051.625
$ws.Range("M91").Value = -1075
$ws.Range("N91").Value = -3859.625
$ws.Range("H125").Value = 715.3333
$ws.Range("J125").Value = 614.8889
$ws.Range("L125").Value = 5534.0001
$ws.Range("N125").Value = -10454.0001
$ws.Range("H127").Value = 1971.1666
$ws.Range("I127").Value = 1971.1666
$ws.Range("K127").Value = 5913.4998
$ws.Range("M127").Value = -953.4997999999996
$ws.Range("H129").Value = 1675.8572
$ws.Range("I129").Value = 840.2857
$ws.Range("K129").Value = 2520.8571
$ws.Range("M129").Value = 2479.1429
$ws.Range("H131").Value = 1065.75
$ws.Range("I131").Value = 1065.75
$ws.Range("K131").Value = 3197.25
$ws.Range("M131").Value = 1842.75
$ws.Range("H132").Value = 2541.258
$ws.Range("I132").Value = 2559.926
$ws.Range("J132").Value = 2415.25
$ws.Range("K132").Value = 7679.778
$ws.Range("L132").Value = 7245.75
$ws.Range("M132").Value = -5149.778
$ws.Range("N132").Value = -12305.75
$ws.Range("H141").Value = 7717.615
$ws.Range("I141").Value = 8281.75
$ws.Range("K141").Value = 24845.25
$ws.Range("M141").Value = -19665.25
# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H45").Value = 2869.4
$ws.Range("I45").Value = 2869.4
$ws.Range("K45").Value = 2869.4
$ws.Range("M45").Value = -2492.4
$ws.Range("H61").Value = 6481.1904
$ws.Range("I61").Value = 5522.143
$ws.Range("J61").Value = 8399.286
$ws.Range("K61").Value = 5522.143
$ws.Range("L61").Value = 8399.286
$ws.Range("M61").Value = -5310.143
$ws.Range("N61").Value = -8823.286
$ws.Range("H74").Value = 1329.6
$ws.Range("I74").Value = 1310.6666
$ws.Range("K74").Value = 1310.6666
$ws.Range("M74").Value = -436.6666
$ws.Range("H77").Value = 1329.6
$ws.Range("I77").Value = 1310.6666
$ws.Range("K77").Value = 6553.333000000001
$ws.Range("M77").Value = -2185.333000000001
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H136").Value = 6481.1904
$ws.Range("I136").Value = 5522.143
$ws.Range("J136").Value = 8399.286
$ws.Range("K136").Value = 16566.429
$ws.Range("L136").Value = 25197.858
$ws.Range("M136").Value = -14016.429
$ws.Range("N136").Value = -30297.858
# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H20").Value = 2288.3
$ws.Range("I20").Value = 2428
$ws.Range("J20").Value = 2148.6
$ws.Range("K20").Value = 2428
$ws.Range("L20").Value = 2148.6
$ws.Range("M20").Value = -2181
$ws.Range("N20").Value = -2642.6
$ws.Range("H35").Value = 63381
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").Value = -104950
# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H4").Value = 2491745.5
$ws.Range("I4").Value = 950182.4399999999
$ws.Range("J4").Value = 7501825
$ws.Range("K4").Value = 2850547.32
$ws.Range("L4").Value = 22505475
$ws.Range("M4").Value = -2850435.32
$ws.Range("N4").Value = -22505699
$ws.Range("H6").Value = 132.3
$ws.Range("I6").Value = 119.22222
$ws.Range("K6").Value = 357.66666
$ws.Range("M6").Value = -244.66666
$ws.Range("H39").Value = 7318
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H58").Value = 500
$ws.Range("I58").Value = 500
$ws.Range("K58").Value = 1500
$ws.Range("M58").Value = -1372
$ws.Range("H75").Value = 461.2857
$ws.Range("J75").Value = 589.8
$ws.Range("L75").Value = 1769.4
$ws.Range("N75").Value = -3765.4
$ws.Range("H78").Value = 461.2857
$ws.Range("J78").Value = 589.8
$ws.Range("L78").Value = 5308.2
$ws.Range("N78").Value = -15292.2
$ws.Range("H80").Value = 6730.8
$ws.Range("J80").Value = 6350
$ws.Range("L80").Value = 19050
$ws.Range("N80").Value = -20922
$ws.Range("H83").Value = 6730.8
$ws.Range("J83").Value = 6350
$ws.Range("L83").Value = 57150
$ws.Range("N83").Value = -66510
$ws.Range("H113").Value = 1854.2
$ws.Range("I113").Value = 1324.8636
$ws.Range("J113").Value = 2750
$ws.Range("K113").Value = 3974.5908
$ws.Range("L113").Value = 8250
$ws.Range("M113").Value = -1804.5908
$ws.Range("N113").Value = -12590
$ws.Range("H114").Value = 1683.3846
$ws.Range("I114").Value = 2338.8
$ws.Range("J114").Value = 1273.75
$ws.Range("K114").Value = 7016.400000000001
$ws.Range("L114").Value = 3821.25
$ws.Range("M114").Value = -3762.400000000001
$ws.Range("N114").Value = -10329.25
$ws.Range("H129").Value = 1436.6
$ws.Range("I129").Value = 1420.875
$ws.Range("J129").Value = 1499.5
$ws.Range("K129").Value = 4262.625
$ws.Range("L129").Value = 4498.5
$ws.Range("M129").Value = 737.375
$ws.Range("N129").Value = -14498.5
$ws.Range("H139").Value = 4385.5293
$ws.Range("I139").Value = 4236.933
$ws.Range("K139").Value = 12710.799
$ws.Range("M139").Value = -7570.798999999999
$ws.Range("H140").Value = 628005.9399999999
$ws.Range("I140").Value = 628005.9399999999
$ws.Range("K140").Value = 1884017.82
$ws.Range("M140").Value = -1878837.82
$ws.Range("H141").Value = 7788.1113
$ws.Range("I141").Value = 7788.1113
$ws.Range("K141").Value = 23364.3339
$ws.Range("M141").Value = -18184.3339
# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H2").Value = 278.7143
$ws.Range("I2").Value = 295.16666
$ws.Range("J2").Value = 180
$ws.Range("K2").Value = 295.16666
$ws.Range("L2").Value = 180
$ws.Range("M2").Value = -182.16666
$ws.Range("N2").Value = -406
$ws.Range("H46").Value = 20343.312
$ws.Range("I46").Value = 5151.25
$ws.Range("K46").Value = 5151.25
$ws.Range("M46").Value = -4995.25
$ws.Range("H80").Value = 2465.875
$ws.Range("I80").Value = 2065.4
$ws.Range("K80").Value = 2065.4
$ws.Range("M80").Value = -1067.4
$ws.Range("H83").Value = 2465.875
$ws.Range("I83").Value = 2065.4
$ws.Range("K83").Value = 10327
$ws.Range("M83").Value = -5335
$ws.Range("H113").Value = 921.8333
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 5266.6
$ws.Range("J122").Value = 5694
$ws.Range("L122").Value = 17082
$ws.Range("N122").Value = -21982
$ws.Range("H126").Value = 3670.6667
$ws.Range("I126").Value = 3006
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 9018
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -6548
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470
# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value = 367
$ws.Range("I22").Value = 300.5
$ws.Range("K22").Value = 300.5
$ws.Range("M22").Value = -5.5
$ws.Range("H27").Value = 367
$ws.Range("I27").Value = 300.5
$ws.Range("K27").Value = 300.5
$ws.Range("M27").Value = -193.5
$ws.Range("H100").Value = 1661
$ws.Range("I100").Value = 1556.5
$ws.Range("K100").Value = 1556.5
$ws.Range("M100").Value = -1015.5
$ws.Range("H132").Value = 3407
$ws.Range("I132").Value = 3224.8333
$ws.Range("K132").Value = 9674.499899999999
$ws.Range("M132").Value = -7144.499899999999
$ws.Range("H136").Value = 6207.077
$ws.Range("I136").Value = 5153.8184
$ws.Range("K136").Value = 15461.4552
$ws.Range("M136").Value = -12911.4552
